$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.75"
$ws.Range("E2").Value = "'1.74%"
$ws.Range("D3").Value = "'41.77"
$ws.Range("E3").Value = "'3.21%"
$ws.Range("D4").Value = "'5.011"
$ws.Range("D5").Value = "'0.07522"
$ws.Range("E5").Value = "'2.49%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.383"
$ws.Range("E6").Value = "'2.02%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.577"
$ws.Range("E7").Value = "'2.21%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9267"
$ws.Range("E8").Value = "'0.48%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.401"
$ws.Range("E9").Value = "'0.05%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1196"
$ws.Range("E10").Value = "'-1.73%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1826"
$ws.Range("E11").Value = "'5.59%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08873"
$ws.Range("E12").Value = "'3.10%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04094"
$ws.Range("E13").Value = "'-4.12%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1046"
$ws.Range("E14").Value = "'-0.66%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001283"
$ws.Range("E15").Value = "'0.28%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006012"
$ws.Range("E16").Value = "'4.01%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.362"
$ws.Range("E17").Value = "'0.68%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3315"
$ws.Range("E18").Value = "'0.85%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'8.079"
$ws.Range("E19").Value = "'4.86%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1392"
$ws.Range("E20").Value = "'0.11%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.3301"
$ws.Range("E21").Value = "'20.13%"
$ws.Range("D22").Value = "'0.04105"
$ws.Range("E22").Value = "'4.42%"
$ws.Range("E23").Value = "'0.38%"
$ws.Range("D24").Value = "'0.003888"
$ws.Range("E24").Value = "'2.95%"
$ws.Range("E25").Value = "'-3.97%"
$ws.Range("D38").Value = "'0.02406"
$ws.Range("E38").Value = "'4.58%"
$ws.Range("D39").Value = "'0.05201"
$ws.Range("E39").Value = "'4.59%"
$ws.Range("E40").Value = "'8.95%"
$ws.Range("D41").Value = "'0.007812"
$ws.Range("E41").Value = "'1.43%"
$ws.Range("D42").Value = "'0.1326"
$ws.Range("E42").Value = "'3.36%"
$ws.Range("D43").Value = "'0.007416"
$ws.Range("E43").Value = "'0.67%"
$ws.Range("D44").Value = "'0.007394"
$ws.Range("E44").Value = "'-5.06%"
$ws.Range("D45").Value = "'0.2950"
$ws.Range("E45").Value = "'-6.96%"
$ws.Range("D46").Value = "'0.00006464"
$ws.Range("E46").Value = "'1.86%"
$ws.Range("D48").Value = "'0.03401"
$ws.Range("E48").Value = "'60.69%"
